$d = $word.ActiveDocument

# Locate the paragraph that begins with "Sodelujete v svetovni ..." (the
# activity-guide intro paragraph) by searching for its unique opening words.
$rng = $d.Content
$found = $rng.Find.Execute("Sodelujete v svetovni", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target paragraph"
}

$p = $rng.Paragraphs(1)

# Insert a throw-away character right at the very start of the paragraph.
# This nudges the paragraph's leading <w:proofErr .../> marker away from the
# absolute start-of-paragraph boundary so that, when we delete+replace the
# whole paragraph text below, it gets swept away together with the rest of
# the old runs instead of being left behind as an orphan.
$p.Range.InsertBefore("X")

$pr = $p.Range
[void]$pr.MoveEnd(1, -1)
$pr.Delete()
$pr.InsertAfter("Sodelujete v svetovni aktivnosti opazovanja in beleženja najšibkejših, s prostim očesom  še vidnih zvezd, kot metode za merjenje svetlobnega onesnaževanja na določenem mestu. Z opazovanjem izbranega Ozvezdje Laboda na nočnem nebu in s primerjavo videnega z zvezdnimi kartami, se lahko ljudje širom sveta podučijo o tem, kako svetila v njihovem kraju prispevajo k svetlobnemu onesnaževanju.  Vaši prispevki v spletno bazo podatkov bodo pomagali dokumentirati nočno nebo, vidno s prostim očesom.")
